$d = $word.ActiveDocument

$bullet = [char]0x2022

# The header/contact-info paragraph ("Rohit Kadam" / "Pune, India ... linkedin...")
# is the document's second paragraph (the first paragraph is a blank spacer line).
$headerRange = $d.Paragraphs(2).Range

# Change 1: the runs "Pune, India <bullet>" and " " (a lone trailing space run)
# are collapsed into a single run "Pune, India <bullet> ".
$old1 = "Pune, India " + $bullet + " "
$headerRange.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# Change 2: the runs " ", "<bullet> +91 7767880235 <bullet>" and " " are collapsed
# into a single run " <bullet> +91 7767880235 <bullet> ".
$old2 = " " + $bullet + " +91 7767880235 " + $bullet + " "
$headerRange.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

Write-Output "Header contact-line runs merged."
